# Changes batch #9.1 Apr 12
#
# The "shapiro" sheet had a spurious "Time bin N" label column repeated
# before the 2nd and 3rd (parameter, statistic, p-value) blocks, plus a
# trailing one after the last block (columns D, H and L). Remove those
# three columns so each block is just (parameter, statistic, p-value).
#
# Deleting from right to left keeps the remaining column letters stable
# while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shapiro")

$ws.Range("L1:L6").EntireColumn.Delete()
$ws.Range("H1:H6").EntireColumn.Delete()
$ws.Range("D1:D6").EntireColumn.Delete()
